$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D currently carries the old "year" date formatting (column-wide
# style + custom width). That formatting needs to move to column C, and
# column D needs to go back to plain/default formatting (it will hold the
# new "solo mujeres" tasa numbers instead of dates). Clear column D's
# formatting first, then restore the header cell's bold/centered look by
# copying the format already used by the other header cells.
$ws.Columns.Item(4).ClearFormats()
$ws.Cells.Item(1, 1).Copy() | Out-Null
$ws.Cells.Item(1, 4).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Swap header labels in C1/D1: C1 becomes "year", D1 becomes "tasa"
$ws.Cells.Item(1, 3).Value = "year"
$ws.Cells.Item(1, 4).Value = "tasa"

# New "solo mujeres" tasa values, keyed by row number
$newTasa = @{
    2  = 27.86
    3  = 80.863
    4  = 94.09399999999999
    5  = 95.502
    6  = 96.10899999999999
    7  = 96.404
    8  = 95.25700000000001
    9  = 92.874
    10 = 86.988
    11 = 61.967
    12 = 22.345
    13 = 18.22
    14 = 71.68091666666666
    15 = 81.70208333333333
    16 = 81.39116666666666
    17 = 82.77375000000001
    18 = 84.29916666666666
    19 = 85.527
    20 = 83.78766666666667
    21 = 78.03416666666666
    22 = 57.98033333333333
    23 = 9.69675
}

for ($row = 2; $row -le 23; $row++) {
    # Column C: now holds the "year" value (date-formatted serial number),
    # which previously lived in column D.
    $cCell = $ws.Cells.Item($row, 3)
    $cCell.Value = 73051
    $cCell.NumberFormat = "yyyy-mm-dd"

    # Column D: now holds the new "solo mujeres" tasa value, with plain
    # (unstyled) formatting -- the same formatting column C used to have.
    $dCell = $ws.Cells.Item($row, 4)
    $dCell.Value = $newTasa[$row]
    $dCell.Style = "Normal"
}

# The custom column width/style formatting moves from column D to column C.
$ws.Columns.Item(3).ColumnWidth = 19.75
